# ADD results from server
# Update computed result values on the "2025", "2030", and "2035" sheets
# (row 2 contains the single data row for each sheet).

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 31.95649908062548
$ws.Range("E2").Value = 290289.2796736782
$ws.Range("I2").Value = 266703.2387598415
$ws.Range("L2").Value = 260645.1078540457
$ws.Range("M2").Value = 117617.5092488
$ws.Range("N2").Value = 71937.21288289552
$ws.Range("O2").Value = 71499.76484539866

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 52084.29526967678
$ws.Range("E2").Value = 269023.3536636419
$ws.Range("I2").Value = 330004.1321038401
$ws.Range("M2").Value = 106735.0832568625
$ws.Range("N2").Value = 36015.51404675592
$ws.Range("O2").Value = 25698.16686044829

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 56692.15729593654
$ws.Range("B2").Value = 28736.2914608956
$ws.Range("E2").Value = 40227.11758927508
$ws.Range("I2").Value = 150567.8529574802
$ws.Range("M2").Value = 58452.06939412496
$ws.Range("N2").Value = 23685.2630177212
$ws.Range("O2").Value = 60683.24652560872
